$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("E2").Value2 = '[Clara%González-Sanguino%NULL%0,              Berta%Ausín%NULL%1,              Miguel Ángel%Castellanos%NULL%1,              Jesús%Saiz%NULL%1,              Aída%López-Gómez%NULL%1,              Carolina%Ugidos%NULL%1,              Manuel%Muñoz%NULL%1]'
$ws.Range("E3").Value2 = '[Cuiyan%Wang%NULL%0,              Riyu%Pan%NULL%0,              Riyu%Pan%NULL%0,              Xiaoyang%Wan%NULL%0,              Yilin%Tan%NULL%0,              Linkang%Xu%NULL%0,              Cyrus S.%Ho%NULL%0,              Roger C.%Ho%NULL%0,              Roger C.%Ho%NULL%0]'
$ws.Range("E4").Value2 = '[Huiyao%Wang%NULL%1,              Qian%Xia%NULL%1,              Zhenzhen%Xiong%NULL%1,              Zhixiong%Li%NULL%1,              Weiyi%Xiang%NULL%1,              Yiwen%Yuan%NULL%1,              Yaya%Liu%NULL%1,              Zhe%Li%NULL%1,              Kenji%Hashimoto%NULL%0,              Kenji%Hashimoto%NULL%0,              Kenji%Hashimoto%NULL%0]'
$ws.Range("E5").Value2 = '[Yenan%Wang%xref no email%0,        Yu%Di%xref no email%0,        Junjie%Ye%xref no email%0,        Wenbin%Wei%xref no email%0]'
$ws.Range("E6").Value2 = '[Yuan%Wu%NULL%1,              Jun%Wang%NULL%1,              Chenggang%Luo%NULL%1,              Sheng%Hu%NULL%1,              Xi%Lin%NULL%1,              Aimee E.%Anderson%NULL%1,              Eduardo%Bruera%NULL%1,              Xiaoxin%Yang%NULL%1,              Shaozhong%Wei%NULL%1,              Yu%Qian%173653835@qq.com%1]'
$ws.Range("E7").Value2 = '[Yingfei%Zhang%NULL%0,              Zheng Feei%Ma%NULL%0,              Zheng Feei%Ma%NULL%0]'
$ws.Range("E9").Value2 = '[Cuiyan%Wang%NULL%0,              Riyu%Pan%NULL%0,              Xiaoyang%Wan%NULL%0,              Yilin%Tan%NULL%0,              Linkang%Xu%NULL%0,              Roger S.%McIntyre%NULL%1,              Faith N.%Choo%NULL%1,              Bach%Tran%NULL%1,              Roger%Ho%NULL%1,              Vijay K.%Sharma%NULL%0,              Cyrus%Ho%NULL%1]'
$ws.Range("E10").Value2 = '[Dana%Tzur Bitan%NULL%2,              Ariella%Grossman-Giron%NULL%1,              Yuval%Bloch%NULL%1,              Yael%Mayer%NULL%1,              Noga%Shiffman%NULL%1,              Shlomo%Mendlovic%NULL%1]'
$ws.Range("E12").Value2 = '[Alison M.%Bacon%ambacon@plymouth.ac.uk%1,              Philip J.%Corr%NULL%2,              Philip J.%Corr%NULL%0]'
$ws.Range("E13").Value2 = '[Shuang-Jiang%Zhou%NULL%4,              Li-Gang%Zhang%NULL%1,              Lei-Lei%Wang%NULL%1,              Zhao-Chang%Guo%NULL%1,              Jing-Qi%Wang%NULL%1,              Jin-Cheng%Chen%NULL%1,              Mei%Liu%NULL%1,              Xi%Chen%NULL%1,              Jing-Xu%Chen%chenjx1110@163.com%1]'
$ws.Range("E14").Value2 = '[Haozheng%Cai%NULL%0,              Baoren%Tu%NULL%0,              Jing%Ma%NULL%0,              Limin%Chen%NULL%0,              Lei%Fu%NULL%0,              Yongfang%Jiang%NULL%0,              Quan%Zhuang%NULL%0]'
$ws.Range("E15").Value2 = '[Naiara%Ozamiz-Etxebarria%xref no email%1,        Maria%Dosil-Santamaria%xref no email%1,        Maitane%Picaza-Gorrochategui%xref no email%1,        Nahia%Idoiaga-Mondragon%xref no email%1]'
$ws.Range("E16").Value2 = '[Lei%Lei%NULL%1,              Xiaoming%Huang%NULL%1,              Shuai%Zhang%NULL%1,              Jinrong%Yang%NULL%1,              Lin%Yang%NULL%0,              Min%Xu%NULL%1]'
$ws.Range("E17").Value2 = '[Shuai%Yuan%NULL%0,              Zhenxin%Liao%NULL%0,              Haojie%Huang%NULL%0,              Boyue%Jiang%NULL%0,              Xueyan%Zhang%NULL%0,              Yingwen%Wang%NULL%0,              Mingyi%Zhao%NULL%0]'
$ws.Range("E18").Value2 = '[Lijun%Kang%NULL%0,              Simeng%Ma%NULL%0,              Min%Chen%NULL%0,              Jun%Yang%NULL%0,              Ying%Wang%NULL%0,              Ruiting%Li%NULL%0,              Lihua%Yao%NULL%0,              Hanping%Bai%NULL%0,              Zhongxiang%Cai%NULL%0,              Bing%Xiang Yang%NULL%0,              Shaohua%Hu%NULL%0,              Kerang%Zhang%NULL%0,              Gaohua%Wang%NULL%0,              Ci%Ma%NULL%0,              Zhongchun%Liu%NULL%0]'
$ws.Range("E19").Value2 = '[Wenjun%Cao%wjcao16@czmc.edu.cn%0,              Ziwei%Fang%NULL%0,              Guoqiang%Hou%NULL%0,              Mei%Han%NULL%0,              Xinrong%Xu%NULL%0,              Jiaxin%Dong%NULL%0,              Jianzhong%Zheng%zhengjj@czmc.edu.cn%0]'
$ws.Range("E20").Value2 = '[Jianbo%Lai%NULL%0,              Simeng%Ma%NULL%0,              Ying%Wang%NULL%0,              Zhongxiang%Cai%NULL%0,              Jianbo%Hu%NULL%0,              Ning%Wei%NULL%0,              Jiang%Wu%NULL%0,              Hui%Du%NULL%0,              Tingting%Chen%NULL%0,              Ruiting%Li%NULL%0,              Huawei%Tan%NULL%0,              Lijun%Kang%NULL%0,              Lihua%Yao%NULL%0,              Manli%Huang%NULL%0,              Huafen%Wang%NULL%0,              Gaohua%Wang%NULL%0,              Zhongchun%Liu%NULL%0,              Shaohua%Hu%NULL%0]'
$ws.Range("E21").Value2 = '[Yuanyuan%Mo%NULL%0,              Lan%Deng%NULL%0,              Lan%Deng%NULL%0,              Liyan%Zhang%NULL%0,              Qiuyan%Lang%NULL%0,              Chunyan%Liao%NULL%0,              Nannan%Wang%NULL%0,              Mingqin%Qin%2026142822@qq.com%0,              Huiqiao%Huang%820325832@qq.com%0]'
$ws.Range("E22").Value2 = '[Sijia%Li%NULL%0,              Yilin%Wang%NULL%1,              Jia%Xue%NULL%1,              Nan%Zhao%NULL%1,              Tingshao%Zhu%NULL%1]'
$ws.Range("E23").Value2 = '[Chen-Yun%Liu%NULL%0,              Yun-zhi%Yang%NULL%1,              Xiao-Ming%Zhang%NULL%1,              Xinying%Xu%NULL%2,              Xinying%Xu%NULL%0,              Qing-Li%Dou%NULL%1,              Wen-Wu%Zhang%NULL%1,              Andy S. K.%Cheng%NULL%1]'
$ws.Range("E24").Value2 = '[Han%Xiao%NULL%0,              Yan%Zhang%NULL%0,              Desheng%Kong%NULL%0,              Shiyue%Li%NULL%0,              Ningxi%Yang%NULL%0]'
$ws.Range("E25").Value2 = '[Derek Gard-Ching%Ying%xref no email%1,        Sunfang%Jiang%xref no email%1,        Hua%Yang%xref no email%1,        Shanzhu%Zhu%xref no email%1]'
$ws.Range("E26").Value2 = '[Paolo%Soraci%paolo.soraci85@gmail.com%2,              Ambra%Ferrari%a.ferrari50@campus.unimib.it%1,              Francesco A.%Abbiati%francesco.abbiati@edu.unito.it%1,              Elena%Del Fante%edelfante@hotmail.it%1,              Rosanna%De Pace%Rosannadepace@libero.it%1,              Antonino%Urso%antonino.urso@tin.it%1,              Mark D.%Griffiths%mark.griffiths@ntu.ac.uk%1]'
